$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# "Marking" row, "Right" column: marks awarded per correct answer (3 -> 5)
$ws.Range("B11").Value = 5

# "Total" row, "Right" column: total marks earned for correct answers (45 -> 75)
$ws.Range("B12").Value = 75

# Correct-marks / total-marks summary label (42/84 -> 75/140)
$ws.Range("E12").Value = "75/140"
